# Sweep_Tables_Examples.xlsx — "added more support for FPGA comm"
#
# The sweep table's second column (B2:B51) held all zeros; populate it
# with a constant sweep value of 10 (one value per row, rows 2-51),
# then leave the selection on that newly-edited block, matching where
# the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill B2:B51 with 10 (was 0 for every row in that block).
$ws.Range("B2:B51").Value = 10

# Reflect the in-progress edit: selection sits on the column just filled.
$ws.Range("B2:B51").Select()
